# Update 'F' column (想去人数 / number of people interested) values
# across the '展览', '演出', and '全部类型' worksheets, as produced by the
# gh-pages data regeneration at commit 456a3b4.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$updates = @{
    3 = 231
    4 = 4822
    5 = 215
    6 = 160
    7 = 121
    10 = 0
    11 = 226
    12 = 0
    14 = 259
    16 = 0
    17 = 152
    18 = 113
    19 = 4035
    20 = 6354
    21 = 39
    23 = 88
    25 = 48
    26 = 3984
    27 = 409
    28 = 48
    29 = 26
    30 = 2587
    32 = 534
    33 = 0
    34 = 0
    36 = 376
    37 = 178
    38 = 0
    41 = 47
    42 = 72
    47 = 78
    48 = 0
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item("演出")
$updates = @{
    2 = 112
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}

$ws = $wb.Worksheets.Item("全部类型")
$updates = @{
    2 = 0
    3 = 231
    4 = 0
    5 = 215
    6 = 160
    7 = 121
    10 = 95
    11 = 762
    12 = 0
    14 = 0
    16 = 187
    17 = 0
    19 = 113
    20 = 4036
    21 = 6354
    24 = 88
    25 = 0
    27 = 3984
    28 = 409
    29 = 48
    30 = 26
    31 = 2587
    33 = 534
    35 = 295
    36 = 0
    37 = 0
    38 = 0
    39 = 10
    40 = 1565
    41 = 972
    42 = 47
    44 = 60
    46 = 482
    49 = 587
}
foreach ($row in $updates.Keys) {
    $ws.Cells.Item([int]$row, 6).Value = $updates[$row]
}
